$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 549, pushing existing rows 549-676 down to 550-677
$ws.Rows.Item(549).Insert()

# Populate the new row 549 with the new data
$ws.Cells.Item(549, 1).Value = 10
$ws.Cells.Item(549, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(549, 3).Value = "La Araucanía"
$ws.Cells.Item(549, 4).Value = 45204
$ws.Cells.Item(549, 5).Value = 9
$ws.Cells.Item(549, 6).Value = 100114014
$ws.Cells.Item(549, 7).Value = "Betarraga"
$ws.Cells.Item(549, 8).Value = "Sin especificar"
$ws.Cells.Item(549, 9).Value = "Primera"
$ws.Cells.Item(549, 10).Value = 250
$ws.Cells.Item(549, 11).Value = 800
$ws.Cells.Item(549, 12).Value = 800
$ws.Cells.Item(549, 13).Value = 800
$ws.Cells.Item(549, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(549, 15).Value = "Región del Maule"
$ws.Cells.Item(549, 16).Value = 67
$ws.Cells.Item(549, 17).Value = 12
$ws.Cells.Item(549, 18).Value = "Hortaliza"
